$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(867).Insert()
$ws.Cells.Item(867, 1).NumberFormat = "@"
$ws.Cells.Item(867, 1).Value = "2026/02/23"
$ws.Cells.Item(867, 1).Style = "Normal"
$ws.Cells.Item(867, 2).Value = "月"
$ws.Cells.Item(867, 3).Value = 16
$ws.Cells.Item(867, 4).Value = 201
